$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Region A (bottom-most first): the "Assumptions" .. "Design pattern" block.
# Paragraph 32 (index at start of script) is the paragraph that currently
# holds only the _GoBack bookmark, right after "Assumptions". Replace its
# content with the full new 9-paragraph sequence, then delete the old
# "Dependencies" / empty / "Constraints" / empty / "High-level design" /
# "System overview" / empty paragraphs that used to follow it (now shifted
# down by the 8 extra paragraphs we just inserted).
# ---------------------------------------------------------------------------
$pBookmark = $d.Paragraphs(32)
Write-Output ("Region A anchor text=[" + $pBookmark.Range.Text + "]")

$regionAXml =
  '<w:p ' + $wns + '>' +
    '<w:r><w:t xml:space="preserve">We don’t assume that the user is familiar with our type of games or has experience with any software. </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>Therefore</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> it needs to be simple (simplicity constraint).</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:r><w:t>We can’t make the assumption that the user has powerful computers (accessibility constraint)</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:pPr><w:pStyle w:val="Heading4"/></w:pPr>' +
    '<w:r><w:t>Dependencies</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:r><w:t>Not sure about this yet</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:pPr><w:pStyle w:val="Heading4"/></w:pPr>' +
    '<w:r><w:t>Constraints</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:r><w:t>Explained above</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:pPr><w:pStyle w:val="Heading3"/></w:pPr>' +
    '<w:r><w:t>High-level design</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:pPr><w:pStyle w:val="Heading4"/></w:pPr>' +
    '<w:r><w:t>System overview</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:pPr><w:pStyle w:val="Heading4"/></w:pPr>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>'

$pBookmark.Range.InsertXML($regionAXml)
Write-Output ("After region A insert, ParaCount=" + $d.Paragraphs.Count)

# The old duplicated paragraphs now sit 8 slots further down (9 new paragraphs
# replaced 1 old one): old #33..#39 are now at #41..#47.
# 41 Dependencies, 42 empty, 43 Constraints, 44 empty(<w:p/>), 45 High-level design,
# 46 System overview, 47 empty (the one that used to precede "Design pattern").
$oldBlockStart = 41
$oldBlockEnd = 47
Write-Output ("Deleting old duplicate paragraphs " + $oldBlockStart + ".." + $oldBlockEnd + ":")
for ($i = $oldBlockStart; $i -le $oldBlockEnd; $i++) {
    Write-Output ("  " + $i + "=[" + $d.Paragraphs($i).Range.Text + "]")
}
$delRange = $d.Range($d.Paragraphs($oldBlockStart).Range.Start, $d.Paragraphs($oldBlockEnd).Range.End)
$delRange.Delete()
Write-Output ("After region A cleanup, ParaCount=" + $d.Paragraphs.Count)

# ---------------------------------------------------------------------------
# Region B: "Public repository for open source, sharing, help from community"
# (paragraph 29) -> split into two runs with amended wording, then two new
# paragraphs inserted after it.
# ---------------------------------------------------------------------------
$pPublic = $d.Paragraphs(29)
Write-Output ("Region B anchor text=[" + $pPublic.Range.Text + "]")
$publicXml =
  '<w:p ' + $wns + '>' +
    '<w:r><w:t>Public repository for open source, help from community</w:t></w:r>' +
    '<w:r><w:t>. we believe community builds better software</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:r><w:t>We believe in sharing</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:r><w:t>We learned from other people’s project, we want our project to be a learning resource for other people</w:t></w:r>' +
  '</w:p>'
$pPublic.Range.InsertXML($publicXml)
Write-Output ("After region B insert, ParaCount=" + $d.Paragraphs.Count)

# ---------------------------------------------------------------------------
# Region C: "MySQL free" (paragraph 27) -> two runs "MySQL " + "because it's
# free" (same combined text, split across two runs).
# ---------------------------------------------------------------------------
$pMysql = $d.Paragraphs(27)
Write-Output ("Region C anchor text=[" + $pMysql.Range.Text + "]")
$mysqlXml =
  '<w:p ' + $wns + '>' +
    '<w:r><w:t xml:space="preserve">MySQL </w:t></w:r>' +
    '<w:r><w:t>because it’s free</w:t></w:r>' +
  '</w:p>'
$pMysql.Range.InsertXML($mysqlXml)
Write-Output ("After region C replace, ParaCount=" + $d.Paragraphs.Count)

# ---------------------------------------------------------------------------
# Region D: after "goal is to make sure any team of programmers can produce
# the wanted result with this document" (paragraph 5), insert a new Heading4
# "Term definition" paragraph and a "Write this at the end" paragraph.
# ---------------------------------------------------------------------------
$pGoal = $d.Paragraphs(5)
Write-Output ("Region D anchor text=[" + $pGoal.Range.Text + "]")
$goalFollowXml =
  '<w:p ' + $wns + '>' +
    '<w:pPr><w:pStyle w:val="Heading4"/></w:pPr>' +
    '<w:r><w:t>Term definition</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:r><w:t>Write this at the end</w:t></w:r>' +
  '</w:p>'
$afterGoal = $pGoal.Range.InsertParagraphAfter()
$d.Paragraphs(6).Range.InsertXML($goalFollowXml)
Write-Output ("After region D insert, ParaCount=" + $d.Paragraphs.Count)

# ---------------------------------------------------------------------------
# Region E: "detail the design of the project" (paragraph 4) -> replace text,
# then insert two new explanatory paragraphs after it.
# ---------------------------------------------------------------------------
$pDetail = $d.Paragraphs(4)
Write-Output ("Region E anchor text=[" + $pDetail.Range.Text + "]")
$detailXml =
  '<w:p ' + $wns + '>' +
    '<w:r><w:t>The purpose of the document is to explain in detail the intricacies of the design choices.</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:r><w:t>Define the different parameters of the project, different rules and constraints.</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wns + '>' +
    '<w:r><w:t>The document is addresses to programmers and will be using a language that is very technical. Will often reference code and assumes that the reader is familiar with software development methods.</w:t></w:r>' +
  '</w:p>'
$pDetail.Range.InsertXML($detailXml)
Write-Output ("After region E replace, ParaCount=" + $d.Paragraphs.Count)

# ---------------------------------------------------------------------------
# Region F: after the second "Introduction" paragraph (Heading3, paragraph 2),
# insert a "Write this at the end" paragraph before "Purpose".
# ---------------------------------------------------------------------------
$pIntro = $d.Paragraphs(2)
Write-Output ("Region F anchor text=[" + $pIntro.Range.Text + "]")
$introFollowXml =
  '<w:p ' + $wns + '>' +
    '<w:r><w:t>Write this at the end</w:t></w:r>' +
  '</w:p>'
$afterIntro = $pIntro.Range.InsertParagraphAfter()
$d.Paragraphs(3).Range.InsertXML($introFollowXml)
Write-Output ("After region F insert, ParaCount=" + $d.Paragraphs.Count)

Write-Output ("FINAL ParaCount=" + $d.Paragraphs.Count)
